$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Round the Ost (Q) / Nord (R) coordinate values to whole numbers for rows 2 and 3
$ws.Range("Q2").Value = 789572
$ws.Range("R2").Value = 7434799
$ws.Range("Q3").Value = 789572
$ws.Range("R3").Value = 7434799

# Clear the Starttid (Z) and Sluttid (AB) cells for rows 2 and 3
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
